$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values must be kept as TEXT (matching the source data's
# inlineStr string cells), so force Text number format before assigning.
$numericCells = "C2","D2","C3","D3","C7","D7","C45","D45","C46","D46","C52","D52","C57","C60","D60","C66","D66","C84","D84","C102","D102","C104","D104","C113","D113"
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Overview / Currency
$ws.Range("C2").Value = "84.42"
$ws.Range("D2").Value = "84.43"

# Row 3 - Overview / Stock Market
$ws.Range("C3").Value = "77691"
$ws.Range("D3").Value = "78675"

# Row 7 - Overview / Inflation Rate
$ws.Range("C7").Value = "6.21"
$ws.Range("D7").Value = "5.49"
$ws.Range("H7").Value = "Oct/24"

# Row 45 - Prices / Inflation Rate
$ws.Range("C45").Value = "6.21"
$ws.Range("D45").Value = "5.49"
$ws.Range("H45").Value = "Oct/24"

# Row 46 - Prices / Consumer Price Index CPI
$ws.Range("C46").Value = "197"
$ws.Range("D46").Value = "194"
$ws.Range("H46").Value = "Oct/24"

# Row 52 - Prices / Food Inflation
$ws.Range("C52").Value = "10.87"
$ws.Range("D52").Value = "9.24"
$ws.Range("H52").Value = "Oct/24"

# Row 57 - Prices / CPI Housing Utilities (D57 unchanged)
$ws.Range("C57").Value = "183"
$ws.Range("H57").Value = "Oct/24"

# Row 58 - Prices / CPI Transportation (only date changes)
$ws.Range("H58").Value = "Oct/24"

# Row 60 - Prices / Inflation Rate MoM
$ws.Range("C60").Value = "1.34"
$ws.Range("D60").Value = "0.6"
$ws.Range("H60").Value = "Oct/24"

# Row 63 - Money / Interbank Rate (only date changes)
$ws.Range("H63").Value = "Nov/24"

# Row 66 - Money / Money Supply M3
$ws.Range("C66").Value = "262159"
$ws.Range("D66").Value = "262940"
$ws.Range("H66").Value = "Oct/24"

# Row 84 - Trade / Crude Oil Production
$ws.Range("C84").Value = "603"
$ws.Range("D84").Value = "605"
$ws.Range("H84").Value = "Jul/24"

# Row 102 - Business / Industrial Production
$ws.Range("C102").Value = "3.1"
$ws.Range("D102").Value = "-0.1"
$ws.Range("H102").Value = "Sep/24"

# Row 104 - Business / Manufacturing Production
$ws.Range("C104").Value = "3.9"
$ws.Range("D104").Value = "1"
$ws.Range("H104").Value = "Sep/24"

# Row 113 - Business / Total Vehicle Sales
$ws.Range("C113").Value = "345107"
$ws.Range("D113").Value = "315689"
$ws.Range("H113").Value = "Oct/24"
